# Applies the "Changes Committed to report and updated data" edit:
#  - adds Solver add-in defined names (left behind by running Solver once)
#  - updates the one-step-ahead harvest input for 2022 (C26) and the
#    dependent LN()/EXP()/MAPE formulas that flow from it
#  - simplifies the per-row MAPE helper formulas (drop the no-op AVERAGE()
#    wrapper around a scalar ABS())
#  - adds a new "best model" regression summary block (tidy()-style
#    term/estimate/std.error/statistic/p.value table) in S33:AA37
#  - tidies up now-unused fill/border formatting on a few ranges
#  - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Solver add-in state (Data > Solver leaves these hidden sheet-scoped
#    defined names behind after being configured/run once).
# ---------------------------------------------------------------------
function Add-HiddenName($name, $refersTo) {
    $n = $ws.Names.Add($name, $refersTo)
    $n.Visible = $false
}

Add-HiddenName "solver_eng" "=1"
Add-HiddenName "solver_neg" "=1"
Add-HiddenName "solver_num" "=0"
Add-HiddenName "solver_opt" "='f_model_one_step_ahead function'!`$I`$4"
Add-HiddenName "solver_typ" "=1"
Add-HiddenName "solver_val" "=0"
Add-HiddenName "solver_ver" "=3"

# ---------------------------------------------------------------------
# 2. Column width hints for the new regression-summary columns.
# ---------------------------------------------------------------------
$ws.Columns("M:M").ColumnWidth = 14.3
$ws.Columns("N:N").ColumnWidth = 14.88
$ws.Columns("P:P").ColumnWidth = 14.88

# ---------------------------------------------------------------------
# 3. Updated 2022 catch input -> ripples through LN/EXP/MAPE columns.
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 18.036332000000002
$ws.Range("X26").Formula = "=LN(C26)"

# ---------------------------------------------------------------------
# 4. Drop the redundant AVERAGE() wrapper on the per-row MAPE helper
#    column (AE22:AE26) -- AVERAGE() of a single scalar is a no-op.
# ---------------------------------------------------------------------
$ws.Range("AE22").Formula = "=(ABS((X22-AA22)/X22))"
$ws.Range("AE23").Formula = "=(ABS((X23-AA23)/X23))"
$ws.Range("AE24").Formula = "=(ABS((X24-AA24)/X24))"
$ws.Range("AE25").Formula = "=(ABS((X25-AA25)/X25))"
$ws.Range("AE26").Formula = "=(ABS((X26-AA26)/X26))"

# ---------------------------------------------------------------------
# 5. Clear the leftover cell formatting on ranges that no longer carry
#    any visible fill/border (their xf entries become unused on save).
# ---------------------------------------------------------------------
foreach ($addr in @("X22:X26", "F33:G36", "F41:K42", "F46:N46")) {
    $r = $ws.Range($addr)
    $r.Borders.LineStyle = -4142
    $r.Interior.Pattern = -4142
}

# J42:K42 were blank style-only placeholder cells -- now completely empty.
$ws.Range("J42:K42").Clear()

# ---------------------------------------------------------------------
# 6. New "best  model" regression summary table (S33:AA37).
# ---------------------------------------------------------------------
$ws.Range("S33").Value = "best  model"

$ws.Range("S34").Value = "model"
$ws.Range("T34").Value = "term"
$ws.Range("U34").Value = "estimate"
$ws.Range("V34").Value = "std.error"
$ws.Range("W34").Value = "statistic"
$ws.Range("X34").Value = "p.value"
$ws.Range("Y34").Value = "sigma"

$ws.Range("S35").Value = "m11"
$ws.Range("T35").Value = "(Intercept)"
$ws.Range("U35").Value = 5.2720784700000003
$ws.Range("V35").Value = 0.67
$ws.Range("W35").Value = 7.8710000000000004
$ws.Range("X35").Value = 0
$ws.Range("Y35").Value = 0.31222536789385502
$ws.Range("Z35").Formula = "=U35+(U36*V27)+(U37*N27)"
$ws.Range("AA35").Formula = "=EXP(Z35)*EXP(0.5*Y35*Y35)"

$ws.Range("S36").Value = "m11"
$ws.Range("T36").Value = "CPUE"
$ws.Range("U36").Value = 0.45925796000000002
$ws.Range("V36").Value = 0.051999999999999998
$ws.Range("W36").Value = 8.8789999999999996
$ws.Range("X36").Value = 0

$ws.Range("S37").Value = "m11"
$ws.Range("T37").Value = "NSEAK_SST_May"
$ws.Range("U37").Value = -0.40041544000000001
$ws.Range("V37").Value = 0.09
$ws.Range("W37").Value = -4.4489999999999998
$ws.Range("X37").Value = 0

# ---------------------------------------------------------------------
# 7. Selection as left by the author.
# ---------------------------------------------------------------------
$ws.Range("P34").Select()
